# Iuvo: switch to English website
#
# The German "Iuvo" results sheets gain a new "Ein-/Auszahlungen" (deposits/
# withdrawals) column, inserted right after the "Endsaldo" column on every
# sheet, filled with 0 for existing data rows. The active sheet/selection
# also moves back to the first sheet ("Tagesergebnisse").

$wb = $excel.ActiveWorkbook

$newHeader = "Ein-/Auszahlungen"

# --- Sheet "Tagesergebnisse": new column F (after E "Endsaldo") ---
$ws1 = $wb.Worksheets.Item("Tagesergebnisse")
$ws1.Columns("F").Insert()
$ws1.Cells.Item(1, 6).Value2 = $newHeader
$ws1.Cells.Item(2, 6).Value2 = 0
$ws1.Columns("F").ColumnWidth = 17.83

# --- Sheet "Monatsergebnisse": new column F (after E "Endsaldo") ---
$ws2 = $wb.Worksheets.Item("Monatsergebnisse")
$ws2.Columns("F").Insert()
$ws2.Cells.Item(1, 6).Value2 = $newHeader
for ($r = 2; $r -le 5; $r++) {
    $ws2.Cells.Item($r, 6).Value2 = 0
}
$ws2.Columns("F").ColumnWidth = 17.83

# --- Sheet "Gesamtergebnis": new column E (after D "Endsaldo") ---
$ws3 = $wb.Worksheets.Item("Gesamtergebnis")
$ws3.Columns("E").Insert()
$ws3.Cells.Item(1, 5).Value2 = $newHeader
for ($r = 2; $r -le 3; $r++) {
    $ws3.Cells.Item($r, 5).Value2 = 0
}
$ws3.Columns("E").ColumnWidth = 17.83

# --- Restore per-sheet selections, then land back on the first tab ---
$null = $ws2.Activate()
$null = $ws2.Range("F9").Select()

$null = $ws3.Activate()
$null = $ws3.Range("E1").Select()

$null = $ws1.Activate()
$null = $ws1.Range("F9").Select()
